$p = $ppt.ActivePresentation

# The presentation currently applies the "Integral" theme's color scheme to
# its (single) slide master -- this is the theme XML part that PowerPoint
# actually renders through. The edit swaps it back to the stock
# "Office Theme" color scheme (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink),
# i.e. it restores the default Office colors in place of Integral's greens.

$master = $p.SlideMaster
$cs = $master.ColorScheme

function HexToRGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2),16)
    $g = [Convert]::ToInt32($hex.Substring(2,2),16)
    $b = [Convert]::ToInt32($hex.Substring(4,2),16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme color scheme, in dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink order (ColorScheme.Colors(1..12)).
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $cs.Colors($i).RGB = HexToRGB($officeThemeColors[$i - 1])
}
